# Auto-generated edit script: update cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Style = "Normal"
    $rng.NumberFormat = "@"
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "69.868.93"
Set-TextCell "E2" "  +0.00%  "
Set-TextCell "D3" "3.531.43"
Set-TextCell "E3" "  +0.75%  "
Set-TextCell "E4" "  +0.08%  "
Set-TextCell "D5" "604.81"
Set-TextCell "E5" "  -0.47%  "
Set-TextCell "D6" "196.37"
Set-TextCell "E6" "  +2.26%  "
Set-TextCell "E7" "  -0.30%  "
Set-TextCell "E8" "  -0.01%  "
Set-TextCell "D9" "0.203"
Set-TextCell "E9" "  -4.82%  "
Set-TextCell "D10" "0.646"
Set-TextCell "E10" "  -2.79%  "
Set-TextCell "D11" "53.45"
Set-TextCell "E11" "  -0.09%  "
Set-TextCell "E12" "  -1.39%  "
Set-TextCell "D13" "9.48"
Set-TextCell "E13" "  -1.57%  "
Set-TextCell "D14" "4.095.80"
Set-TextCell "E14" "  +0.79%  "
Set-TextCell "D15" "598.29"
Set-TextCell "E15" "  -3.64%  "
Set-TextCell "B16" "WrappedBTC"
Set-TextCell "C16" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell "D16" "69.997.05"
Set-TextCell "E16" "  +0.08%  "
Set-TextCell "B17" "Uniswap"
Set-TextCell "C17" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell "D17" "12.74"
Set-TextCell "E17" "  +0.17%  "
Set-TextCell "D18" "19.01"
Set-TextCell "E18" "  +0.67%  "
Set-TextCell "D19" "3.530.56"
Set-TextCell "E19" "  +0.90%  "
Set-TextCell "E20" "  +1.83%  "
Set-TextCell "D21" "0.985"
Set-TextCell "E21" "  -0.70%  "
Set-TextCell "D22" "17.94"
Set-TextCell "E22" "  +1.45%  "
Set-TextCell "D23" "103.27"
Set-TextCell "E23" "  -2.48%  "
Set-TextCell "D24" "5.16"
Set-TextCell "E24" "  +3.01%  "
Set-TextCell "D25" "4.63"
Set-TextCell "E25" "  -0.33%  "
Set-TextCell "E26" "  +0.82%  "
Set-TextCell "D27" "10.80"
Set-TextCell "E27" "  -1.82%  "
Set-TextCell "D28" "9.55"
Set-TextCell "E28" "  -2.98%  "
Set-TextCell "D29" "33.30"
Set-TextCell "E29" "  -2.58%  "
Set-TextCell "E30" "  -0.41%  "
Set-TextCell "D31" "4.23"
Set-TextCell "E31" "  +1.44%  "
Set-TextCell "E32" "  -2.06%  "
Set-TextCell "E33" "  -0.07%  "
Set-TextCell "D34" "63.50"
Set-TextCell "E34" "  -1.07%  "
Set-TextCell "B35" "Maker"
Set-TextCell "C35" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D35" "3.777.51"
Set-TextCell "E35" "  +1.49%  "
Set-TextCell "B36" "Fetch.AI"
Set-TextCell "C36" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D36" "3.16"
Set-TextCell "E36" "  +2.01%  "
Set-TextCell "D37" "0.0₃0813"
Set-TextCell "E37" "  +2.40%  "
Set-TextCell "E38" "  +0.08%  "
Set-TextCell "B39" "Bittensor"
Set-TextCell "C39" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D39" "506.71"
Set-TextCell "E39" "  -4.23%  "
Set-TextCell "B40" "TheGraph"
Set-TextCell "C40" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell "D40" "0.391"
Set-TextCell "E40" "  -0.01%  "
Set-TextCell "D41" "3.58"
Set-TextCell "E41" "  -0.08%  "
Set-TextCell "D42" "36.47"
Set-TextCell "E42" "  -0.90%  "
Set-TextCell "D43" "0.133"
Set-TextCell "E43" "  -3.09%  "
Set-TextCell "E44" "  -2.90%  "
Set-TextCell "E45" "  -0.94%  "
Set-TextCell "E46" "  -1.57%  "
Set-TextCell "E47" "  -2.53%  "
Set-TextCell "D48" "1.01"
Set-TextCell "E48" "  +0.25%  "
Set-TextCell "D49" "8.48"
Set-TextCell "E49" "  -3.08%  "
Set-TextCell "B50" "Mantle"
Set-TextCell "C50" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D50" "1.34"
Set-TextCell "E50" "  +3.02%  "
Set-TextCell "B51" "FLOKI"
Set-TextCell "C51" "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextCell "D51" "0.000248"
Set-TextCell "E51" "  +5.39%  "
